# Update the public EPEX Spot / Gaz / CO2 price export with the newest day.
#
# "Prix Spot" sheet: append a new date column (BR) "22-aug" with its 24
#   hourly prices (rows 2-25), copying the formatting of the previous day
#   column (BQ) so the new column matches the existing table style.
# "Gaz" and "CO2" sheets: append a new row (67) for 2025-08-20 with its
#   "Last Price" value.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Prix Spot": new column BR = 22-aug
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Prix Spot")

# Copy style/borders/font from the last existing day column (BQ) onto the
# new one (BR) so the appended column looks like the rest of the table.
$ws.Range("BQ1:BQ25").Copy($ws.Range("BR1:BR25"))

$ws.Range("BR1").Value = "22-aug"

# Hourly prices for rows 2..25, in row order (00-01 .. 23-24).
$hourlyValues = @(
    66.38,
    51.66,
    48.66,
    41.88,
    31.51,
    36.41,
    71.03,
    70.88,
    75.09999999999999,
    49.67,
    15.57,
    5.11,
    19.02,
    6.83,
    5.37,
    10.07,
    23.6,
    26.45,
    73.48,
    84.56,
    104.53,
    103.49,
    101.59,
    94.37
)

for ($i = 0; $i -lt $hourlyValues.Count; $i++) {
    $row = $i + 2
    $ws.Range("BR$row").Value = $hourlyValues[$i]
}

# ---------------------------------------------------------------------
# Sheet "Gaz": new row 67 = 2025-08-20 / 30.85
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Gaz")

$cell = $ws.Range("A67")
$cell.NumberFormat = "@"          # force text so the date-like string isn't
$cell.Value = "2025-08-20"        # auto-converted into a date serial number
$cell.ClearFormats()              # drop the temporary "@" format again so the
                                   # cell is left with the sheet's default style

$ws.Range("B67").Value = 30.85

# ---------------------------------------------------------------------
# Sheet "CO2": new row 67 = 2025-08-20 / 70.3
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CO2")

$cell = $ws.Range("A67")
$cell.NumberFormat = "@"
$cell.Value = "2025-08-20"
$cell.ClearFormats()

$ws.Range("B67").Value = 70.3
